{"js": "// The commit duplicates the last table row (OCR text repeated three times,\n// \"\u0414\u0430 \u044d\u0442\u043e \u0437\u043d\u0430\u0447\u0435\u043d\u0438\u0435 \u043a\u043b\u044e\u0447\u0430 \u0414\u0430 \u044d\u0442\u043e \u0437\u043d\u0430\u0447\u0435\u043d\u0438\u0435 \u043a\u043b\u044e\u0447\u0430 \u0414\u0430 \u044d\u0442\u043e \u0437\u043d\u0430\u0447\u0435\u043d\u0438\u0435 \u043a\u043b\u044e\u0447\u0430, \u043e\u0442 02,09.2020\")\n// two more times at the bottom of the table, so the recognized text ends up\n// appearing on three consecutive rows instead of one.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  return \"no tables found\";\n}\n\nconst table = tables.items[0];\ntable.rows.load(\"items\");\nawait context.sync();\n\nconst rows = table.rows.items;\nconst lastRow = rows[rows.length - 1];\nlastRow.cells.load(\"items\");\nawait context.sync();\n\nconst cells = lastRow.cells.items;\nfor (const cell of cells) {\n  cell.body.load(\"text\");\n}\nawait context.sync();\n\n// Capture the last row's text, cell by cell, so the duplicated rows match\n// exactly (including the repeated OCR text in the 4th column).\nconst lastRowValues = cells.map((cell) => cell.body.text.replace(/\\r?\\n$/, \"\"));\n\n// Append two duplicate rows (identical content to the current last row) to\n// the end of the table.\ntable.addRows(\"End\", 2, [lastRowValues, lastRowValues]);\nawait context.sync();\n\nreturn \"done\";\n", "ps1": "# The commit duplicates the last table row (OCR text repeated three times,\n# \"\u0414\u0430 \u044d\u0442\u043e \u0437\u043d\u0430\u0447\u0435\u043d\u0438\u0435 \u043a\u043b\u044e\u0447\u0430 \u0414\u0430 \u044d\u0442\u043e \u0437\u043d\u0430\u0447\u0435\u043d\u0438\u0435 \u043a\u043b\u044e\u0447\u0430 \u0414\u0430 \u044d\u0442\u043e \u0437\u043d\u0430\u0447\u0435\u043d\u0438\u0435 \u043a\u043b\u044e\u0447\u0430, \u043e\u0442 02,09.2020\")\n# two more times at the bottom of the table, so the recognized text ends up\n# appearing on three consecutive rows instead of one.\n\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\n$colCount = $tbl.Columns.Count\n$lastRow = $tbl.Rows.Item($tbl.Rows.Count)\n\n# Capture the current last row's cell text (column by column) so the\n# duplicated rows match exactly, including the repeated OCR text.\n$cellTexts = @()\nfor ($c = 1; $c -le $colCount; $c++) {\n    $cellTexts += $lastRow.Cells.Item($c).Range.Text\n}\n\n# Insert two duplicate rows directly above the existing last row, each one\n# a copy of that row's content.\nfor ($i = 0; $i -lt 2; $i++) {\n    $anchorRow = $tbl.Rows.Item($tbl.Rows.Count)\n    $newRow = $tbl.Rows.Add($anchorRow)\n    for ($c = 1; $c -le $colCount; $c++) {\n        $newRow.Cells.Item($c).Range.Text = $cellTexts[$c - 1]\n    }\n}\n"}
